# Weekly update: insert a new "Poroto granado" price record for
# Vega Central Mapocho de Santiago ahead of the existing history,
# shifting the older rows down by one (the data set keeps growing
# week over week, newest entry goes on top of this sub-range).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 275..331 down to 276..332 to make room for the new record.
$ws.Rows("275:275").Insert()

# Populate the newly inserted row 275 with this week's data point.
$ws.Cells.Item(275, 1).Value  = 9
$ws.Cells.Item(275, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(275, 3).Value  = "Metropolitana"
$ws.Cells.Item(275, 4).Value  = 44932
$ws.Cells.Item(275, 5).Value  = 13
$ws.Cells.Item(275, 6).Value  = 100112030
$ws.Cells.Item(275, 7).Value  = "Poroto granado"
$ws.Cells.Item(275, 8).Value  = "Sin especificar"
$ws.Cells.Item(275, 9).Value  = "Primera"
$ws.Cells.Item(275, 10).Value = 70
$ws.Cells.Item(275, 11).Value = 40000
$ws.Cells.Item(275, 12).Value = 43000
$ws.Cells.Item(275, 13).Value = 41500
$ws.Cells.Item(275, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(275, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(275, 16).Value = 1660
$ws.Cells.Item(275, 17).Value = 25
$ws.Cells.Item(275, 18).Value = "Hortaliza"
